$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: "time taken (minutes)" ---
# Copy header formatting (bold style) from O1 into P1
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "time taken (minutes)"
$ws.Rows("1").RowHeight = 24.45

# Time-taken values for existing rows 2-4
$ws.Range("P2").Value = 12
$ws.Range("P3").Value = 12
$ws.Range("P4").Value = 28

# --- New row 5: R18_imagenet_v2 ---
# Copy formatting from row 4 down into row 5 (keeps the s=1 / s=4 style pattern)
$ws.Range("A4:O4").Copy()
$ws.Range("A5:O5").PasteSpecial(-4122)
$ws.Rows("5").RowHeight = 57.45

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "R18_imagenet_v2"
$ws.Range("C5").Value = "resnet18-pretrained, 512-128 linear, relu, 128-5 linear, softmax"
$ws.Range("D5").Value = "pretrained  model weights, fc layer alone trained"
$ws.Range("E5").Value = "Adam"
$ws.Range("F5").Value = "NLLLoss"
$ws.Range("G5").Value = "RandomResizedCrop(448)"
$ws.Range("H5").Value = "Dataset mean, std normalise"
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 0.0004
$ws.Range("K5").Value = "One cycle lr"
$ws.Range("L5").Value = "StratifiedKFold (5 splits)"
$ws.Range("M5").Value = "[0.95946, 0.74856, 0.68570]"
$ws.Range("N5").Value = "[0.69729, 0.74860, 0.77215]"
$ws.Range("O5").Value = "[0.75425, 0.66836, 0.63191]"
$ws.Range("P5").Value = 20.5

$ws.Range("A6").Select() | Out-Null
